# Bold the "BAB 1 :" heading and the "PERMULAAN" subheading, and relocate
# the stray "_GoBack" bookmark from the end of the first section to the
# blank paragraph that follows "PERMULAAN".

$d = $word.ActiveDocument

# Paragraph 2: "B" + "AB 1 :" -- bold the whole paragraph (runs + the
# paragraph mark itself, i.e. pPr/rPr) in one shot.
$babPara = $d.Paragraphs.Item(2)
$babPara.Range.Font.Bold = 1

# Paragraph 3: "PERMULAAN" -- same treatment.
$permulaanPara = $d.Paragraphs.Item(3)
$permulaanPara.Range.Font.Bold = 1

# The blank paragraph right after "PERMULAAN" gets the "_GoBack" bookmark.
# Word keeps only one bookmark per name, so adding it here automatically
# removes it from wherever it used to live (the last paragraph of the
# first section, right after the sectPr).
$blankPara = $d.Paragraphs.Item(4)
$d.Bookmarks.Add("_GoBack", $blankPara.Range)
